$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("set commands")

$ws.Range("A529").Value = "set shared profiles url-filtering Outbound-URL credential-enforcement alert real-time-detection"
$ws.Range("A530").Value = "set shared profiles url-filtering Outbound-URL alert real-time-detection"
$ws.Range("A531").Value = "set shared profiles url-filtering Alert-Only-URL credential-enforcement mode ip-user"
$ws.Range("A532").Value = "set shared profiles url-filtering Alert-Only-URL credential-enforcement log-severity medium"
$ws.Range("A533").Value = "set shared profiles url-filtering Alert-Only-URL credential-enforcement alert [ Block Allow abortion abused-drugs adult alcohol-and-tobacco auctions business-and-economy command-and-control computer-and-internet-info content-delivery-networks copyright-infringement cryptocurrency dating dynamic-dns educational-institutions entertainment-and-arts extremism financial-services gambling games government grayware hacking health-and-medicine high-risk home-and-garden hunting-and-fishing insufficient-content internet-communications-and-telephony internet-portals job-search legal low-risk malware medium-risk military motor-vehicles music newly-registered-domain news not-resolved nudity online-storage-and-backup parked peer-to-peer personal-sites-and-blogs philosophy-and-political-advocacy phishing private-ip-addresses proxy-avoidance-and-anonymizers questionable real-estate recreation-and-hobbies reference-and-research religion search-engines sex-education shareware-and-freeware shopping social-networking society sports stock-advice-and-tools streaming-media swimsuits-and-intimate-apparel training-and-tools translation travel unknown weapons web-advertisements web-based-email web-hosting ]"
$ws.Range("A534").Value = "set shared profiles url-filtering Alert-Only-URL alert [ Block Allow abortion abused-drugs adult alcohol-and-tobacco auctions business-and-economy command-and-control computer-and-internet-info content-delivery-networks copyright-infringement cryptocurrency dating dynamic-dns educational-institutions entertainment-and-arts extremism financial-services gambling games government grayware hacking health-and-medicine high-risk home-and-garden hunting-and-fishing insufficient-content internet-communications-and-telephony internet-portals job-search legal low-risk malware medium-risk military motor-vehicles music newly-registered-domain news not-resolved nudity online-storage-and-backup parked peer-to-peer personal-sites-and-blogs philosophy-and-political-advocacy phishing private-ip-addresses proxy-avoidance-and-anonymizers questionable real-estate recreation-and-hobbies reference-and-research religion search-engines sex-education shareware-and-freeware shopping social-networking society sports stock-advice-and-tools streaming-media swimsuits-and-intimate-apparel training-and-tools translation travel unknown weapons web-advertisements web-based-email web-hosting ]"
$ws.Range("A535").Value = "set shared profiles url-filtering Alert-Only-URL mlav-engine-urlbased-enabled `"`"Phishing Detection`"`" mlav-policy-action alert"
$ws.Range("A536").Value = "set shared profiles url-filtering Alert-Only-URL mlav-engine-urlbased-enabled `"`"Javascript Exploit Detection`"`" mlav-policy-action alert"
$ws.Range("A537").Value = "set shared profiles url-filtering Alert-Only-URL credential-enforcement alert real-time-detection"
$ws.Range("A538").Value = "set shared profiles url-filtering Alert-Only-URL alert real-time-detection"
$ws.Range("A539").Value = "set shared profiles url-filtering Exception-URL credential-enforcement mode ip-user"
$ws.Range("A540").Value = "set shared profiles url-filtering Exception-URL credential-enforcement log-severity high"
$ws.Range("A541").Value = "set shared profiles url-filtering Exception-URL credential-enforcement block [ Block Allow abortion abused-drugs adult alcohol-and-tobacco auctions business-and-economy command-and-control computer-and-internet-info content-delivery-networks copyright-infringement cryptocurrency dating dynamic-dns educational-institutions entertainment-and-arts extremism financial-services gambling games government grayware hacking health-and-medicine high-risk home-and-garden hunting-and-fishing insufficient-content internet-communications-and-telephony internet-portals job-search legal low-risk malware medium-risk military motor-vehicles music newly-registered-domain news not-resolved nudity online-storage-and-backup parked peer-to-peer personal-sites-and-blogs philosophy-and-political-advocacy phishing private-ip-addresses proxy-avoidance-and-anonymizers questionable real-estate recreation-and-hobbies reference-and-research religion search-engines sex-education shareware-and-freeware shopping social-networking society sports stock-advice-and-tools streaming-media swimsuits-and-intimate-apparel training-and-tools translation travel unknown weapons web-advertisements web-based-email web-hosting ]"
$ws.Range("A542").Value = "set shared profiles url-filtering Exception-URL log-http-hdr-user-agent yes"
$ws.Range("A543").Value = "set shared profiles url-filtering Exception-URL log-http-hdr-referer yes"
$ws.Range("A544").Value = "set shared profiles url-filtering Exception-URL log-http-hdr-xff yes"
$ws.Range("A545").Value = "set shared profiles url-filtering Exception-URL alert [ Allow abortion abused-drugs adult alcohol-and-tobacco auctions business-and-economy computer-and-internet-info content-delivery-networks copyright-infringement cryptocurrency dating dynamic-dns educational-institutions entertainment-and-arts extremism financial-services gambling games government hacking health-and-medicine high-risk home-and-garden hunting-and-fishing insufficient-content internet-communications-and-telephony internet-portals job-search legal low-risk medium-risk military motor-vehicles music newly-registered-domain news not-resolved nudity online-storage-and-backup parked peer-to-peer personal-sites-and-blogs philosophy-and-political-advocacy private-ip-addresses proxy-avoidance-and-anonymizers questionable real-estate recreation-and-hobbies reference-and-research religion search-engines sex-education shareware-and-freeware shopping social-networking society sports stock-advice-and-tools streaming-media swimsuits-and-intimate-apparel training-and-tools translation travel unknown weapons web-advertisements web-based-email web-hosting ]"
$ws.Range("A546").Value = "set shared profiles url-filtering Exception-URL block [ Block command-and-control grayware malware phishing ]"
$ws.Range("A547").Value = "set shared profiles url-filtering Exception-URL mlav-engine-urlbased-enabled `"`"Phishing Detection`"`" mlav-policy-action block"
$ws.Range("A548").Value = "set shared profiles url-filtering Exception-URL mlav-engine-urlbased-enabled `"`"Javascript Exploit Detection`"`" mlav-policy-action block"
